$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep the Price/Volume columns as literal text (they already are, as
# non-numeric inline strings such as "43.255.81") so values like "1.00"
# or "6.20" are not silently coerced into numbers and lose their trailing
# zeros / formatting.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '43.202.49'
$ws.Range("E2").Value = '  +0.13%  '

$ws.Range("D3").Value = '2.323.35'
$ws.Range("E3").Value = '  +0.73%  '

$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").Value = '302.54'
$ws.Range("E5").Value = '  +0.44%  '

$ws.Range("D6").Value = '99.97'
$ws.Range("E6").Value = '  +1.32%  '

$ws.Range("D7").Value = '0.507'
$ws.Range("E7").Value = '  +0.59%  '

$ws.Range("E8").Value = '  +0.06%  '

$ws.Range("D9").Value = '0.520'
$ws.Range("E9").Value = '  +2.50%  '

$ws.Range("D10").Value = '36.66'
$ws.Range("E10").Value = '  +7.75%  '

$ws.Range("E11").Value = '  -0.32%  '

$ws.Range("E12").Value = '  +0.55%  '

$ws.Range("D13").Value = '17.91'
$ws.Range("E13").Value = '  +1.18%  '

$ws.Range("E14").Value = '  +2.77%  '

$ws.Range("D15").Value = '2.686.21'
$ws.Range("E15").Value = '  +0.65%  '

$ws.Range("D16").Value = '2.366.93'

$ws.Range("D17").Value = '0.801'
$ws.Range("E17").Value = '  -0.99%  '

$ws.Range("D18").Value = '43.145.24'
$ws.Range("E18").Value = '  +0.29%  '

$ws.Range("D19").Value = '12.57'
$ws.Range("E19").Value = '  +4.45%  '

$ws.Range("B20").Value = 'Uniswap'
$ws.Range("C20").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D20").Value = '6.20'
$ws.Range("E20").Value = '  +1.57%  '

$ws.Range("B21").Value = 'ShibaInu'
$ws.Range("C21").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D21").Value = '0.0₃0909'
$ws.Range("E21").Value = '  +0.27%  '

$ws.Range("E22").Value = '  +0.56%  '

$ws.Range("D23").Value = '236.99'
$ws.Range("E23").Value = '  -0.20%  '

$ws.Range("D24").Value = '2.19'
$ws.Range("E24").Value = '  +5.60%  '

$ws.Range("D25").Value = '2.46'
$ws.Range("E25").Value = '  +0.09%  '

$ws.Range("D26").Value = '0.999'
$ws.Range("E26").Value = '  -0.04%  '

$ws.Range("D27").Value = '25.42'
$ws.Range("E27").Value = '  +2.96%  '

$ws.Range("D28").Value = '168.53'
$ws.Range("E28").Value = '  +1.66%  '

$ws.Range("D29").Value = '34.59'
$ws.Range("E29").Value = '  +2.00%  '

$ws.Range("D30").Value = '9.18'
$ws.Range("E30").Value = '  +0.05%  '

$ws.Range("D31").Value = '2.05'
$ws.Range("E31").Value = '  -10.40%  '

$ws.Range("D32").Value = '5.20'
$ws.Range("E32").Value = '  +3.70%  '

$ws.Range("D33").Value = '1.00'
$ws.Range("E33").Value = '  -0.04%  '

$ws.Range("B34").Value = 'RenderToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D34").Value = '4.71'
$ws.Range("E34").Value = '  +3.45%  '

$ws.Range("B35").Value = 'Celestia'
$ws.Range("C35").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D35").Value = '17.64'
$ws.Range("E35").Value = '  +3.02%  '

$ws.Range("E36").Value = '  -0.87%  '

$ws.Range("D37").Value = '0.0697'
$ws.Range("E37").Value = '  -0.24%  '

$ws.Range("D38").Value = '0.102'
$ws.Range("E38").Value = '  +0.62%  '

$ws.Range("E40").Value = '  -1.43%  '

$ws.Range("E41").Value = '  +0.68%  '

$ws.Range("D42").Value = '2.001.16'
$ws.Range("E42").Value = '  +0.21%  '

$ws.Range("D43").Value = '0.0291'
$ws.Range("E43").Value = '  +1.68%  '

$ws.Range("E44").Value = '  -4.30%  '

$ws.Range("D45").Value = '10.12'
$ws.Range("E45").Value = '  +0.99%  '

$ws.Range("D46").Value = '17.99'
$ws.Range("E46").Value = '  +0.28%  '

$ws.Range("D47").Value = '2.91'
$ws.Range("E47").Value = '  +1.63%  '

$ws.Range("D48").Value = '55.39'
$ws.Range("E48").Value = '  +2.17%  '

$ws.Range("B49").Value = 'RocketPoolETH'
$ws.Range("C49").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D49").Value = '2.551.12'
$ws.Range("E49").Value = '  +0.75%  '

$ws.Range("B50").Value = 'Stacks'
$ws.Range("C50").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D50").Value = '1.55'
$ws.Range("E50").Value = '  +2.14%  '

$ws.Range("D51").Value = '71.76'
$ws.Range("E51").Value = '  +2.08%  '
